# Lightning2.0RemainingTasks.xlsx - "fix a bug in basic example"
#
# The "basic example" items (Bold text / Italic text / Fix underline-strikethrough
# text bugs) are removed from the Pre-RC task list, along with the now-irrelevant
# "NativeAOT support" row. Two other items ("No more isempty" and
# "Renderer-independent blend mode") together with "Fast texture atlas rendering &
# general optimisation" move down into the already-completed "Pre-Final Tasks"
# section. The "last updated" note at the top is refreshed and restyled to match
# the bold section headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash the formatting of the little blank marker cell C19 (sits on the blank
# row right below the "Pre-Final Tasks" header) so we can restore it onto the
# new blank separator row later, instead of letting it drift away with the
# row deletes below.
$ws.Range("C19").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Remove the four "basic example" related rows (10-13): "Fast texture atlas
# rendering & general optimisation", "Bold text", "Italic text", "Fix
# underline / strikethrough text bugs" - and the three rows (16-18):
# "No more isempty", "Renderer-independent blend mode", "NativeAOT support".
# Deleting bottom-to-top keeps the row numbers for the not-yet-deleted rows
# stable while we work.
$ws.Rows(18).Delete()
$ws.Rows(17).Delete()
$ws.Rows(16).Delete()
$ws.Rows(13).Delete()
$ws.Rows(12).Delete()
$ws.Rows(11).Delete()
$ws.Rows(10).Delete()

# Refresh the "last updated" note and make it look like the bold section
# headers (still keeps its date-ish number format), bump the row height to
# match the other header rows.
$ws.Range("A2").Value = "Updated March 21, 2023"
$ws.Range("A2").Font.Bold = $true
$ws.Rows(2).RowHeight = 15

# "2.0 Examples" (now row 13, used to be row 20 with no "Completed" marker)
# picks up a "Yes" just like its neighbours in the Pre-RC list.
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("B13").Value = "Yes"

# The little blank marker cell drifted to C12 along with the row shifts above
# - clear it out, it doesn't belong there any more.
$ws.Range("C12").Clear()

# Append the three relocated items to the end of the "Pre-Final Tasks"
# section, matching the plain (no "Completed" marker) style used by their new
# neighbours.
$ws.Range("A20").Copy() | Out-Null
$ws.Range("A23:A25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("A23").Value = "No more isempty"
$ws.Range("A24").Value = "Renderer-independent blend mode"
$ws.Range("A25").Value = "Fast texture atlas rendering & general optimisation"

# Restore the blank marker cell's formatting onto the new blank separator row
# (row 19, between the two section headers) and clean up the scratch cell.
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").Clear()

# Move the active selection up to the refreshed date note.
$ws.Range("A2").Select()
